$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master Biaya")

# Insert a new column before column A (shifts B..J right); this becomes
# the new "Use" flag column.
$ws.Columns.Item(1).Insert()

# Header for the new column
$ws.Range("A1").Value = "Use"

# Existing data row (now row 2 after the shift): update Jabatan (col I)
# from "Grade C" to "Grade D".
$ws.Range("I2").Value = "Grade D"

# Duplicate row 2 into a new row 3, keeping all values but with a
# different Jabatan (col H -> "Senior Manager").
$ws.Range("A2:J2").Copy()
$ws.Range("A3").PasteSpecial(-4104)

# Preserve the percentage number format on F3 (copy already brings it,
# but make sure explicitly too).
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("A3").Value = "Yes"
$ws.Range("H3").Value = "Senior Manager"
$ws.Range("A2").Value = "No"

# Add three more blank rows (4-6), each with the percentage style
# carried in column F (matching the data-entry rows below).
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F6").PasteSpecial(-4122)

# Extend the in-sheet data validation lists from row 2 only to rows 2-6.
$ws.Range("B2").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C2:C6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D2:D6").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E2:E6").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F2:F6").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("G2:G6").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("J2:J6").PasteSpecial(-4122)

$ws.Range("B2:B6").Validation.Delete()
$ws.Range("B2:B6").Validation.Add(3, 1, 1, "Luar Negeri,Dalam Negeri,Luar Kota,Dalam Kota")
$ws.Range("B2:B6").Validation.InputMessage = ""
$ws.Range("B2:B6").Validation.ErrorMessage = ""

$ws.Range("C2:C6").Validation.Delete()
$ws.Range("C2:C6").Validation.Add(3, 1, 1, "Transportasi,Uang Makan,Hotel")

$ws.Range("D2:D6").Validation.Delete()
$ws.Range("D2:D6").Validation.Add(3, 1, 1, "1500,1000")

$ws.Range("E2:E6").Validation.Delete()
$ws.Range("E2:E6").Validation.Add(3, 1, 1, "Iya,Tidak")

$ws.Range("F2:F6").Validation.Delete()
$ws.Range("F2:F6").Validation.Add(3, 1, 1, "100%,75%,0")

$ws.Range("G2:G6").Validation.Delete()
$ws.Range("G2:G6").Validation.Add(3, 1, 1, "Aktif,Tidak Aktif")

$ws.Range("J2:J6").Validation.Delete()
$ws.Range("J2:J6").Validation.Add(3, 1, 1, "> 16,>= 12 <= 16,< 12")

# New validation list for the "Use" column.
$ws.Range("A2:A6").Validation.Add(3, 1, 1, "Yes,No")

# Extend the two cross-sheet (x14) data validations from a single row
# to rows 2-6 as well.
$ws.Range("H2:H6").Validation.Delete()
$ws.Range("H2:H6").Validation.Add(3, 1, 1, "=Sheet2!$A$1:$A$9")

$ws.Range("I2:I6").Validation.Delete()
$ws.Range("I2:I6").Validation.Add(3, 1, 1, "=Sheet2!$B$4:$D$4")

# Column A gets no custom width (default), matching the source edit --
# nothing further required there since it was never set.

# Move the active selection to I5, matching the saved view state.
$ws.Range("I5").Select()

Write-Host "done"
